# Auto-generated edit script: updates cached price/profit values on several
# sheets to match the refreshed scheduled-runner data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 221.76471
$ws.Range("I11").Value = 221.76471
$ws.Range("K11").Value = 221.76471
$ws.Range("M11").Value = -81.76471000000001
$ws.Range("H32").Value = 2377.348
$ws.Range("J32").Value = 2805.889
$ws.Range("L32").Value = 2805.889
$ws.Range("N32").Value = -3457.889
$ws.Range("H40").Value = 1166.6666
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H41").Value = 455.30768
$ws.Range("I41").Value = 759.4
$ws.Range("K41").Value = 759.4
$ws.Range("M41").Value = -319.4
$ws.Range("H51").Value = 4758.875
$ws.Range("I51").Value = 5814.2
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 5814.2
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -5330.2
$ws.Range("N51").Value = -3968
$ws.Range("H76").Value = 10560.5
$ws.Range("J76").Value = 10560.5
$ws.Range("L76").Value = 10560.5
$ws.Range("N76").Value = -11190.5
$ws.Range("H79").Value = 10560.5
$ws.Range("J79").Value = 10560.5
$ws.Range("L79").Value = 10560.5
$ws.Range("N79").Value = -12744.5
$ws.Range("H97").Value = 9742
$ws.Range("J97").Value = 9742
$ws.Range("L97").Value = 29226
$ws.Range("N97").Value = -30218
$ws.Range("H111").Value = 706.381
$ws.Range("I111").Value = 286.2857
$ws.Range("K111").Value = 858.8571000000001
$ws.Range("M111").Value = 2208.1429
$ws.Range("H116").Value = 6710.4287
$ws.Range("J116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("N116").Value = -11884
$ws.Range("H132").Value = 11639.55
$ws.Range("I132").Value = 7423.75
$ws.Range("K132").Value = 22271.25
$ws.Range("M132").Value = -19741.25
$ws.Range("H137").Value = 17815.834
$ws.Range("I137").Value = 5270.909
$ws.Range("K137").Value = 15812.727
$ws.Range("M137").Value = -13262.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5287.161
$ws.Range("I32").Value = 2979.149
$ws.Range("J32").Value = 12518.934
$ws.Range("K32").Value = 2979.149
$ws.Range("L32").Value = 12518.934
$ws.Range("M32").Value = -2692.149
$ws.Range("N32").Value = -13092.934
$ws.Range("H61").Value = 782623.6
$ws.Range("I61").Value = 5489.55
$ws.Range("K61").Value = 5489.55
$ws.Range("M61").Value = -5277.55
$ws.Range("H63").Value = 2859.3333
$ws.Range("I63").Value = 2689
$ws.Range("J63").Value = 3200
$ws.Range("K63").Value = 2689
$ws.Range("L63").Value = 3200
$ws.Range("M63").Value = -2003
$ws.Range("N63").Value = -4572
$ws.Range("H66").Value = 2859.3333
$ws.Range("I66").Value = 2689
$ws.Range("J66").Value = 3200
$ws.Range("K66").Value = 13445
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = -10013
$ws.Range("N66").Value = -22864
$ws.Range("H122").Value = 3626
$ws.Range("I122").Value = 3168.3
$ws.Range("J122").Value = 4541.4
$ws.Range("K122").Value = 9504.900000000001
$ws.Range("L122").Value = 13624.2
$ws.Range("M122").Value = -7054.900000000001
$ws.Range("N122").Value = -18524.2
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H136").Value = 782623.6
$ws.Range("I136").Value = 5489.55
$ws.Range("K136").Value = 16468.65
$ws.Range("M136").Value = -13918.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3888.5293
$ws.Range("J86").Value = 4998.778
$ws.Range("L86").Value = 4998.778
$ws.Range("N86").Value = -7244.778
$ws.Range("H89").Value = 3888.5293
$ws.Range("J89").Value = 4998.778
$ws.Range("L89").Value = 24993.89
$ws.Range("N89").Value = -36225.89
$ws.Range("H94").Value = 1712.1482
$ws.Range("I94").Value = 1378.4706
$ws.Range("J94").Value = 2279.4
$ws.Range("K94").Value = 1378.4706
$ws.Range("L94").Value = 2279.4
$ws.Range("M94").Value = -927.4706000000001
$ws.Range("N94").Value = -3181.4
$ws.Range("H99").Value = 14390.0625
$ws.Range("J99").Value = 1637.6
$ws.Range("L99").Value = 1637.6
$ws.Range("N99").Value = -4633.6
$ws.Range("H134").Value = 11107.275
$ws.Range("I134").Value = 6661.4062
$ws.Range("J134").Value = 28890.75
$ws.Range("K134").Value = 19984.2186
$ws.Range("L134").Value = 86672.25
$ws.Range("M134").Value = -17449.2186
$ws.Range("N134").Value = -91742.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 55104.523
$ws.Range("I31").Value = 114107.555
$ws.Range("J31").Value = 17174
$ws.Range("K31").Value = 114107.555
$ws.Range("L31").Value = 17174
$ws.Range("M31").Value = -113812.555
$ws.Range("N31").Value = -17764
$ws.Range("H34").Value = 55104.523
$ws.Range("I34").Value = 114107.555
$ws.Range("J34").Value = 17174
$ws.Range("K34").Value = 114107.555
$ws.Range("L34").Value = 17174
$ws.Range("M34").Value = -113905.555
$ws.Range("N34").Value = -17578
$ws.Range("H105").Value = 13911.583
$ws.Range("I105").Value = 15696
$ws.Range("J105").Value = 11413.4
$ws.Range("K105").Value = 15696
$ws.Range("L105").Value = 11413.4
$ws.Range("M105").Value = -13949
$ws.Range("N105").Value = -14907.4
$ws.Range("H122").Value = 7235.8
$ws.Range("J122").Value = 19253.334
$ws.Range("L122").Value = 57760.00199999999
$ws.Range("N122").Value = -62660.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1714.29
$ws.Range("I131").Value = 860.2222
$ws.Range("J131").Value = 1798.7583
$ws.Range("K131").Value = 2580.6666
$ws.Range("L131").Value = 5396.2749
$ws.Range("M131").Value = 2459.3334
$ws.Range("N131").Value = -15476.2749

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2474
$ws.Range("I46").Value = 805.8333
$ws.Range("K46").Value = 805.8333
$ws.Range("M46").Value = -617.8333
$ws.Range("H55").Value = 596.8570999999999
$ws.Range("I55").Value = 496.42856
$ws.Range("J55").Value = 697.2857
$ws.Range("K55").Value = 496.42856
$ws.Range("L55").Value = 697.2857
$ws.Range("M55").Value = -323.42856
$ws.Range("N55").Value = -1043.2857
$ws.Range("H68").Value = 16863.637
$ws.Range("I68").Value = 13550.1
$ws.Range("J68").Value = 49999
$ws.Range("K68").Value = 13550.1
$ws.Range("L68").Value = 49999
$ws.Range("M68").Value = -12801.1
$ws.Range("N68").Value = -51497
$ws.Range("H71").Value = 16863.637
$ws.Range("I71").Value = 13550.1
$ws.Range("J71").Value = 49999
$ws.Range("K71").Value = 67750.5
$ws.Range("L71").Value = 249995
$ws.Range("M71").Value = -64006.5
$ws.Range("N71").Value = -257483
$ws.Range("H82").Value = 3690.3684
$ws.Range("I82").Value = 1094.8
$ws.Range("J82").Value = 6574.3335
$ws.Range("K82").Value = 1094.8
$ws.Range("L82").Value = 6574.3335
$ws.Range("M82").Value = -733.8
$ws.Range("N82").Value = -7296.3335
$ws.Range("H85").Value = 3690.3684
$ws.Range("I85").Value = 1094.8
$ws.Range("J85").Value = 6574.3335
$ws.Range("K85").Value = 1094.8
$ws.Range("L85").Value = 6574.3335
$ws.Range("M85").Value = 153.2
$ws.Range("N85").Value = -9070.333500000001
$ws.Range("H100").Value = 5173.125
$ws.Range("I100").Value = 2672.9092
$ws.Range("K100").Value = 2672.9092
$ws.Range("M100").Value = -2131.9092
$ws.Range("H132").Value = 1047576.5
$ws.Range("I132").Value = 3303.238
$ws.Range("J132").Value = 2801955.8
$ws.Range("K132").Value = 9909.714
$ws.Range("L132").Value = 8405867.399999999
$ws.Range("M132").Value = -7379.714
$ws.Range("N132").Value = -8410927.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12021.923
$ws.Range("I62").Value = 11599.6
$ws.Range("K62").Value = 11599.6
$ws.Range("M62").Value = -10975.6
$ws.Range("H65").Value = 12021.923
$ws.Range("I65").Value = 11599.6
$ws.Range("K65").Value = 57998
$ws.Range("M65").Value = -54878
$ws.Range("H126").Value = 7764.25
$ws.Range("J126").Value = 11088.909
$ws.Range("L126").Value = 33266.727
$ws.Range("N126").Value = -38206.727
